$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Column F/G (Max/Min Voltage) unit: [%] -> [p.u.]
    $ws.Range("F7").Value = "[p.u.]"
    $ws.Range("G7").Value = "[p.u.]"

    # Description wording tweak: "node" -> "it"
    $ws.Range("P5").Value = "Which package this belongs to"
    $ws.Range("K5").Value = "Year where it is commissioned (1.1.xxxx)"
    $ws.Range("L5").Value = "Year where it is decommissioned (31.12.xxxx)"

    # Database-behavior details for Data Package / Data Source columns
    $ws.Range("P6").Value = "Scenario-dependent"
    $ws.Range("Q6").Value = "Scenario-dependent"
}
